# update database and change read_price algorithm
#
# The balance sheet adds a new reporting period ("12 ماهه منتهی به 1401/12",
# published "1402-02-23"). Every period column shifts one slot to the left
# (D<-E, E<-F, F<-G, G<-H) and the freed-up column H receives the freshly
# scraped figures for the new period. A handful of subtotal rows were also
# recomputed by the updated read_price algorithm, so their shifted G value
# is corrected afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Shift the two header rows (financial period / publish date) left,
#    then drop in the new period's labels in column H.
# ---------------------------------------------------------------------
$ws.Range("D8:G8").Value = $ws.Range("E8:H8").Value()
$ws.Range("D9:G9").Value = $ws.Range("E9:H9").Value()

$ws.Cells.Item(8, 8).Value = "12 ماهه منتهی به 1401/12"
$ws.Cells.Item(9, 8).Value = "1402-02-23 (2)"
$ws.Cells.Item(9, 7).Value = "1402-02-23 (10)"

# ---------------------------------------------------------------------
# 2) Shift every data row's figures one column to the left (D:G <- E:H).
# ---------------------------------------------------------------------
for ($r = 12; $r -le 58; $r++) {
    $shifted = $ws.Range("E$r`:H$r").Value()
    if ($shifted -ne $null) {
        $ws.Range("D$r`:G$r").Value = $shifted
    }
}

# ---------------------------------------------------------------------
# 3) Drop in the freshly reported figures for the new period (column H).
# ---------------------------------------------------------------------
$newPeriod = @{
    12 = 7429894;   13 = 14446292;  14 = 1831783;   15 = 19651905;
    16 = 46320627;  17 = 24458;     18 = 89704959;  19 = 215422;
    20 = 20490;     21 = 0;         22 = 13911365;  23 = 5039;
    24 = "-";       25 = 0;         26 = 14325174;  27 = 104030133;
    29 = 7783732;   30 = "-";       31 = 27677703;  32 = 4165407;
    33 = 33087;     34 = 8436789;   35 = 4150111;   36 = 0;
    37 = 52246829;  38 = 0;         39 = "-";       40 = 15000000;
    41 = 583432;    42 = 15583432;  43 = 67830261;  45 = 16900000;
    46 = 0;         47 = 0;         48 = 0;         49 = 0;
    50 = 1690000;   51 = 20000;     52 = "-";       53 = 0;
    54 = "-";       55 = 0;         56 = 17589872;  57 = 36199872;
    58 = 104030133
}
foreach ($r in $newPeriod.Keys) {
    $ws.Cells.Item($r, 8).Value = $newPeriod[$r]
}

# ---------------------------------------------------------------------
# 4) A few subtotal rows recompute slightly differently than a pure
#    shift under the updated read_price algorithm - correct column G.
# ---------------------------------------------------------------------
$gOverrides = @{
    29 = 4962467;
    31 = 8606726;
    37 = 20499397;
    43 = 25811308;
    56 = 11082602;
    57 = 19761058
}
foreach ($r in $gOverrides.Keys) {
    $ws.Cells.Item($r, 7).Value = $gOverrides[$r]
}
